$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new entry for "Duy" on row 7 (shifting the previously blank row 7 content)
$ws.Range("A7").Value = "Duy"
$ws.Range("B7").Value = 44672
$ws.Range("B7").NumberFormat = "d-mmm"
$ws.Range("C7").Value = 2

# Update the selected cell to match the author's final selection
$ws.Range("D10").Select()
